# Applies the row-level corrections from the "error solve ifrs list" commit:
# rows 2-6 get corrected figures (and lose their "U" column + some trailing
# ratio columns), rows 7-9 lose all their financial figures entirely,
# keeping only the identifying A/B/C columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2301
$ws.Range("E2").Value = 78
$ws.Range("F2").Value = 78
$ws.Range("G2").Value = -495
$ws.Range("H2").Value = -534
$ws.Range("I2").Value = -507
$ws.Range("J2").Value = -28
$ws.Range("K2").Value = 10438
$ws.Range("L2").Value = 6661
$ws.Range("M2").Value = 3777
$ws.Range("N2").Value = 3763
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 3530
$ws.Range("Q2").Value = -966
$ws.Range("R2").Value = -328
$ws.Range("S2").Value = 1511
$ws.Range("T2").Value = 63
$ws.Range("V2").Value = 3713
$ws.Range("W2").Value = 3.38
$ws.Range("X2").Value = -23.21
$ws.Range("Y2").Value = -12.62
$ws.Range("Z2").Value = -5.23
$ws.Range("AA2").Value = 176.36
$ws.Range("AB2").Value = 8.48
$ws.Range("AC2").Value = -718
$ws.Range("AD2").Value = -2.53
$ws.Range("AE2").Value = 5434
$ws.Range("AF2").Value = 0.33
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 60314092
$ws.Range("U2").Value = $null

# Row 3
$ws.Range("D3").Value = 2394
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 359
$ws.Range("H3").Value = 315
$ws.Range("I3").Value = 313
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 10204
$ws.Range("L3").Value = 5992
$ws.Range("M3").Value = 4212
$ws.Range("N3").Value = 4204
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 3530
$ws.Range("Q3").Value = 428
$ws.Range("R3").Value = 170
$ws.Range("S3").Value = -431
$ws.Range("T3").Value = 7
$ws.Range("V3").Value = 3904
$ws.Range("W3").Value = 4.17
$ws.Range("X3").Value = 13.14
$ws.Range("Y3").Value = 7.87
$ws.Range("Z3").Value = 3.04
$ws.Range("AA3").Value = 142.27
$ws.Range("AB3").Value = 20.81
$ws.Range("AC3").Value = 444
$ws.Range("AD3").Value = 5.1
$ws.Range("AE3").Value = 6072
$ws.Range("AF3").Value = 0.37
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 60314092
$ws.Range("U3").Value = $null

# Row 4
$ws.Range("D4").Value = 2631
$ws.Range("E4").Value = 287
$ws.Range("F4").Value = 287
$ws.Range("G4").Value = 390
$ws.Range("H4").Value = 364
$ws.Range("I4").Value = 363
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 13596
$ws.Range("L4").Value = 9017
$ws.Range("M4").Value = 4579
$ws.Range("N4").Value = 4540
$ws.Range("O4").Value = 39
$ws.Range("P4").Value = 3530
$ws.Range("Q4").Value = -1496
$ws.Range("R4").Value = 169
$ws.Range("S4").Value = 2158
$ws.Range("T4").Value = 21
$ws.Range("V4").Value = 5761
$ws.Range("W4").Value = 10.92
$ws.Range("X4").Value = 13.85
$ws.Range("Y4").Value = 8.300000000000001
$ws.Range("Z4").Value = 3.05
$ws.Range("AA4").Value = 196.94
$ws.Range("AB4").Value = 29.71
$ws.Range("AC4").Value = 514
$ws.Range("AD4").Value = 5.69
$ws.Range("AE4").Value = 6431
$ws.Range("AF4").Value = 0.45
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 60314092
$ws.Range("U4").Value = $null

# Row 5
$ws.Range("D5").Value = 3342
$ws.Range("E5").Value = 428
$ws.Range("F5").Value = 428
$ws.Range("G5").Value = 421
$ws.Range("H5").Value = 363
$ws.Range("I5").Value = 360
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 18889
$ws.Range("L5").Value = 13986
$ws.Range("M5").Value = 4903
$ws.Range("N5").Value = 4844
$ws.Range("O5").Value = 59
$ws.Range("P5").Value = 3530
$ws.Range("Q5").Value = -2813
$ws.Range("R5").Value = -12
$ws.Range("S5").Value = 3413
$ws.Range("T5").Value = 22
$ws.Range("V5").Value = 8220
$ws.Range("W5").Value = 12.8
$ws.Range("X5").Value = 10.87
$ws.Range("Y5").Value = 7.66
$ws.Range("Z5").Value = 2.21
$ws.Range("AA5").Value = 285.23
$ws.Range("AB5").Value = 38.91
$ws.Range("AC5").Value = 509
$ws.Range("AD5").Value = 7.65
$ws.Range("AE5").Value = 6862
$ws.Range("AF5").Value = 0.57
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 60314092
$ws.Range("U5").Value = $null
$ws.Range("AG5").Value = $null
$ws.Range("AH5").Value = $null

# Row 6
$ws.Range("D6").Value = 3587
$ws.Range("E6").Value = 371
$ws.Range("F6").Value = 371
$ws.Range("G6").Value = 443
$ws.Range("H6").Value = 344
$ws.Range("I6").Value = 338
$ws.Range("K6").Value = 17507
$ws.Range("L6").Value = 12239
$ws.Range("M6").Value = 5269
$ws.Range("N6").Value = 5196
$ws.Range("P6").Value = 3530
$ws.Range("Q6").Value = 2189
$ws.Range("R6").Value = 263
$ws.Range("S6").Value = -1250
$ws.Range("T6").Value = 30
$ws.Range("V6").Value = 7345
$ws.Range("W6").Value = 10.35
$ws.Range("X6").Value = 9.58
$ws.Range("Y6").Value = 6.73
$ws.Range("Z6").Value = 1.86
$ws.Range("AA6").Value = 232.28
$ws.Range("AB6").Value = 49.27
$ws.Range("AC6").Value = 478
$ws.Range("AD6").Value = 6.19
$ws.Range("AE6").Value = 7360
$ws.Range("AF6").Value = 0.4
$ws.Range("AJ6").Value = 60314092
$ws.Range("U6").Value = $null
$ws.Range("AG6").Value = $null
$ws.Range("AH6").Value = $null
$ws.Range("AI6").Value = $null

# Rows 7-9: every figure column (D through AJ) is removed, only A/B/C remain
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
